$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Work Hours (B) and Work Description (C) for rows 8-14, which were
# previously blank (only the Date in column A was present).
#
# The Work Description strings are written in the same order the author
# originally added them to the shared-string table (C10, C11, C9, C8, C14)
# so new entries land at the same shared-string indices as the target file.
$ws.Range("C10").Value = "Waiting for info from Riemer & Anders"
$ws.Range("C11").Value = "Meeting with Riemer & Anders, License, Started working on the Interactive Tree View"
$ws.Range("C9").Value = "Fixed hour proposition & Emailed Riemer and Anders about meeting."
$ws.Range("C8").Value = "Worked on READ ME and Instructions."
$ws.Range("C12").Value = "Weekend"
$ws.Range("C13").Value = "Weekend"
$ws.Range("C14").Value = "Worked on tab selection and disabling editing when not in the editor tab. "

$ws.Range("B8").Value = 8
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 8
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 8

# Update the view: scroll down a bit and move the selection to C14 (matches
# the author continuing to fill the timesheet further down the sheet).
$ws.Range("C14").Select()
$excel.ActiveWindow.ScrollRow = 4

# The Total: formula in B38 (=SUM(B2:B36)) will recalc automatically to 58
# now that rows 8-14 carry Work Hours values.
